# SITEWORK-257: updated knolltextiles_info and _approvals data files.
#
# - Removes the stray "Amazon" leather-grade row (row 278: Amazon / " " / W / null)
# - Removes the 13 trailing "blank"/"Z" placeholder rows (rows 299-311)
# - Excel auto-shrinks the Table1 ref/autoFilter, sheet dimension and the
#   shared-strings table as a consequence of the row deletions.
# - Tidies up the now-duplicated conditional-formatting rules that were left
#   over (there were two copies of the same 3 "header" + 2 "pending" rules,
#   plus a now-empty A262:E279 block) and re-points the remaining rules at
#   their correct ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the stray "Amazon" row and the trailing placeholder rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(278).Delete()
$ws.Range("A298:A310").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 2. Clean up conditional formatting.
#    Before: 8 rule-groups (40 rules collapsed into 20 cfRule COM items):
#      1-3   A1:E1 D92:D93 D94:E261 A2:B261 A280:E1048576 D2:E91  (Ltd Yes/No/Yes)
#      4-5   A1:E1 D92:D93 D94:E261 A2:C261 A280:E1048576 D2:E91  (Pending x2)
#      6-8   A262:E279                                            (Ltd Yes/No/Yes)
#      9-10  A262:E279                                            (Pending x2)
#      11-13 G261                                                 (Ltd Yes/No/Yes)
#      14-15 G261                                                 (Pending x2)
#      16-18 G262:G297                                            (Ltd Yes/No/Yes)
#      19-20 G262:G297                                            (Pending x2)
#    After: groups 1-10 (the duplicated header rules + the now-obsolete
#    A262:E279 block) are removed; the remaining 4 groups are re-pointed so
#    the "header" rules cover the whole sheet body (A262:E1048576 /
#    G262:G296 folded in) and the G261/G262:G297 rules swap ranges.
# ---------------------------------------------------------------------------

$fcs = $ws.Cells.FormatConditions

# Grab references to the 10 rules we are going to keep *before* deleting
# anything in front of them (indices shift otherwise).
$keepLtdYes1 = $fcs.Item(11)   # dxfId 9 - "Limited Yes"
$keepNo1     = $fcs.Item(12)   # dxfId 8 - "No"
$keepYes1    = $fcs.Item(13)   # dxfId 7 - "Yes"
$keepPend1a  = $fcs.Item(14)   # dxfId 6 - "Pending"
$keepPend1b  = $fcs.Item(15)   # dxfId 5 - "Pending"
$keepLtdYes2 = $fcs.Item(16)   # dxfId 4 - "Limited Yes"
$keepNo2     = $fcs.Item(17)   # dxfId 3 - "No"
$keepYes2    = $fcs.Item(18)   # dxfId 2 - "Yes"
$keepPend2a  = $fcs.Item(19)   # dxfId 1 - "Pending"
$keepPend2b  = $fcs.Item(20)   # dxfId 0 - "Pending"

# Remember the exact look of each rule (font/fill) so it can be re-created
# on the extra ranges that ModifyAppliesToRange can't reach in one shot
# (this COM layer only keeps the first area of a multi-area range).
function Copy-FcLook($src, $dst) {
    $dst.Font.Color = $src.Font.Color
    $dst.Interior.Pattern = $src.Interior.Pattern
    if ($src.Interior.Pattern -ne -4142) {
        $dst.Interior.Color = $src.Interior.Color
    }
}

# Drop the 10 duplicated / obsolete rules (the two A262:E279 blocks and the
# two header blocks that get folded into keepLtdYes1/keepNo1/... below).
for ($i = 1; $i -le 10; $i++) {
    $fcs.Item(1).Delete()
}

# ---- Group C / D: G262:G297 -> G261 --------------------------------------
$g261 = $ws.Range("G261")
$keepLtdYes2.ModifyAppliesToRange($g261)
$keepNo2.ModifyAppliesToRange($g261)
$keepYes2.ModifyAppliesToRange($g261)
$keepPend2a.ModifyAppliesToRange($g261)
$keepPend2b.ModifyAppliesToRange($g261)

# ---- Group A / B: G261 -> big combined range ------------------------------
$areasA = @("A1:E1", "D92:D93", "D94:E261", "A2:B261", "D2:E91", "G262:G296", "A262:E1048576")
$areasB = @("A1:E1", "D92:D93", "D94:E261", "A2:C261", "D2:E91", "G262:G296", "A262:E1048576")

# Point the surviving rule objects at the last (largest / most significant)
# area, then replicate the same look across the remaining areas with fresh
# rules so every cell that should be covered still is.
$lastA = $areasA[$areasA.Length - 1]
$keepLtdYes1.ModifyAppliesToRange($ws.Range($lastA))
$keepNo1.ModifyAppliesToRange($ws.Range($lastA))
$keepYes1.ModifyAppliesToRange($ws.Range($lastA))

$lastB = $areasB[$areasB.Length - 1]
$keepPend1a.ModifyAppliesToRange($ws.Range($lastB))
$keepPend1b.ModifyAppliesToRange($ws.Range($lastB))

for ($i = 0; $i -lt $areasA.Length - 1; $i++) {
    $rng = $ws.Range($areasA[$i])

    $c = $rng.FormatConditions.Add(1, 3, '="Limited Yes"')
    Copy-FcLook $keepLtdYes1 $c

    $c = $rng.FormatConditions.Add(1, 3, '="No"')
    Copy-FcLook $keepNo1 $c

    $c = $rng.FormatConditions.Add(1, 3, '="Yes"')
    Copy-FcLook $keepYes1 $c
}

for ($i = 0; $i -lt $areasB.Length - 1; $i++) {
    $rng = $ws.Range($areasB[$i])

    $c = $rng.FormatConditions.Add(1, 3, '="Pending"')
    Copy-FcLook $keepPend1a $c

    $c = $rng.FormatConditions.Add(1, 3, '="Pending"')
    Copy-FcLook $keepPend1b $c
}
